$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Settings")
$assets = $wb.Worksheets.Item("Assets")

# --- New rows 13-18: create new strings in the exact order required so that
# the shared-strings table gets rebuilt with the same ordering as the target.
$ws.Range("A13").Value = "O365AppID"
$ws.Range("A14").Value = "O365TenantID"
$ws.Range("B13").Value = "Shared_O365ApplicationID"
$ws.Range("B14").Value = "Shared_O365TenantID"
$ws.Range("B15").Value = "Shared_O365ApplicationSecret"
$ws.Range("A15").Value = "O365ApplicationSecret"

# Existing rows whose values change to a new string (template paths switched
# from .docx to .xlsx)
$ws.Range("B6").Value = "C:\Users\56382C\Documents\Hiring Manager Template.xlsx"
$ws.Range("B7").Value = "C:\Users\56382C\Documents\Welcome Email Template.xlsx"

$ws.Range("A16").Value = "ExternalTransferCurrency"
$ws.Range("A17").Value = "ExternalTransferFrequency"
$ws.Range("B16").Value = "USD"
$ws.Range("B17").Value = "Annual"

$ws.Range("A18").Value = "CompCodeOverride"

# DEV value in B3, and the new queue name in B2 (must be assigned after the
# above, and B3 before B2, to preserve shared-string ordering)
$ws.Range("B3").Value = "DEV"
$ws.Range("B2").Value = "P004_SP003_090_NHC_WD_Performer_Queue"

# Numeric value with an explicit left-aligned style (creates new cellXf)
$ws.Range("B18").Value = 8810
$ws.Range("B18").HorizontalAlignment = -4131

# --- Selections ---
$assets.Activate()
$assets.Range("A2:B4").Select()

$ws.Activate()
$ws.Range("B3").Select()
